$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 183.33333
$ws.Range("I4").Value = 183.33333
$ws.Range("K4").Value = 183.33333
$ws.Range("M4").Value = -69.33332999999999
$ws.Range("H33").Value = 809.5909
$ws.Range("I33").Value = 413
$ws.Range("J33").Value = 1503.625
$ws.Range("K33").Value = 413
$ws.Range("L33").Value = 1503.625
$ws.Range("M33").Value = -184
$ws.Range("N33").Value = -1961.625
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825
$ws.Range("H41").Value = 123.27273
$ws.Range("I41").Value = 45.5
$ws.Range("J41").Value = 140.55556
$ws.Range("K41").Value = 45.5
$ws.Range("L41").Value = 140.55556
$ws.Range("M41").Value = 394.5
$ws.Range("N41").Value = -1020.55556
$ws.Range("H62").Value = 6392.625
$ws.Range("I62").Value = 5831.5557
$ws.Range("J62").Value = 7114
$ws.Range("K62").Value = 5831.5557
$ws.Range("L62").Value = 7114
$ws.Range("M62").Value = -5207.5557
$ws.Range("N62").Value = -8362
$ws.Range("H65").Value = 6392.625
$ws.Range("I65").Value = 5831.5557
$ws.Range("J65").Value = 7114
$ws.Range("K65").Value = 29157.7785
$ws.Range("L65").Value = 35570
$ws.Range("M65").Value = -26037.7785
$ws.Range("N65").Value = -41810
$ws.Range("H69").Value = 6013
$ws.Range("I69").Value = 6013
$ws.Range("K69").Value = 18039
$ws.Range("M69").Value = -17165
$ws.Range("H72").Value = 6013
$ws.Range("I72").Value = 6013
$ws.Range("K72").Value = 54117
$ws.Range("M72").Value = -49749
$ws.Range("H74").Value = 5349.857
$ws.Range("I74").Value = 3724.5
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3724.5
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -2788.5
$ws.Range("N74").Value = -7872
$ws.Range("H76").Value = 11328.733
$ws.Range("I76").Value = 13212.9
$ws.Range("K76").Value = 13212.9
$ws.Range("M76").Value = -12897.9
$ws.Range("H77").Value = 5349.857
$ws.Range("I77").Value = 3724.5
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 18622.5
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -13942.5
$ws.Range("N77").Value = -39360
$ws.Range("H79").Value = 11328.733
$ws.Range("I79").Value = 13212.9
$ws.Range("K79").Value = 13212.9
$ws.Range("M79").Value = -12120.9
$ws.Range("H86").Value = 9599.799999999999
$ws.Range("I86").Value = 15000
$ws.Range("K86").Value = 15000
$ws.Range("M86").Value = -13877
$ws.Range("H89").Value = 9599.799999999999
$ws.Range("I89").Value = 15000
$ws.Range("K89").Value = 75000
$ws.Range("M89").Value = -69384
$ws.Range("H116").Value = 15076.25
$ws.Range("J116").Value = 21153
$ws.Range("L116").Value = 21153
$ws.Range("N116").Value = -28037
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1291.0465
$ws.Range("I132").Value = 987.9
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 2963.7
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -433.6999999999998
$ws.Range("N132").Value = -21059
$ws.Range("H135").Value = 8066202.5
$ws.Range("I135").Value = 1126.5454
$ws.Range("K135").Value = 10138.9086
$ws.Range("M135").Value = -7603.908599999999
$ws.Range("H138").Value = 5612.205
$ws.Range("I138").Value = 2488.5
$ws.Range("J138").Value = 6418.3228
$ws.Range("K138").Value = 7465.5
$ws.Range("L138").Value = 19254.9684
$ws.Range("M138").Value = -2325.5
$ws.Range("N138").Value = -29534.9684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H32").Value = 14934874
$ws.Range("I32").Value = 18873162
$ws.Range("K32").Value = 18873162
$ws.Range("M32").Value = -18872875
$ws.Range("H35").Value = 1327.5
$ws.Range("I35").Value = 330
$ws.Range("J35").Value = 2325
$ws.Range("K35").Value = 330
$ws.Range("L35").Value = 2325
$ws.Range("M35").Value = 76
$ws.Range("N35").Value = -3137
$ws.Range("H45").Value = 2024.5
$ws.Range("I45").Value = 1616.6
$ws.Range("J45").Value = 2432.4
$ws.Range("K45").Value = 1616.6
$ws.Range("L45").Value = 2432.4
$ws.Range("M45").Value = -1239.6
$ws.Range("N45").Value = -3186.4
$ws.Range("H74").Value = 143019890
$ws.Range("I74").Value = 143019890
$ws.Range("K74").Value = 143019890
$ws.Range("M74").Value = -143019016
$ws.Range("H77").Value = 143019890
$ws.Range("I77").Value = 143019890
$ws.Range("K77").Value = 715099450
$ws.Range("M77").Value = -715095082
$ws.Range("H97").Value = 602.34784
$ws.Range("I97").Value = 375.94736
$ws.Range("J97").Value = 1677.75
$ws.Range("K97").Value = 375.94736
$ws.Range("L97").Value = 1677.75
$ws.Range("M97").Value = 120.05264
$ws.Range("N97").Value = -2669.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2602.8333
$ws.Range("I25").Value = 2633
$ws.Range("J25").Value = 2542.5
$ws.Range("K25").Value = 2633
$ws.Range("L25").Value = 2542.5
$ws.Range("M25").Value = -2398
$ws.Range("N25").Value = -3012.5
$ws.Range("H94").Value = 3018
$ws.Range("I94").Value = 1069.8
$ws.Range("K94").Value = 1069.8
$ws.Range("M94").Value = -618.8
$ws.Range("H99").Value = 3591.3333
$ws.Range("I99").Value = 1906.8636
$ws.Range("J99").Value = 5771.2354
$ws.Range("K99").Value = 1906.8636
$ws.Range("L99").Value = 5771.2354
$ws.Range("M99").Value = -408.8635999999999
$ws.Range("N99").Value = -8767.2354
$ws.Range("H102").Value = 9977
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 7741.3335
$ws.Range("I107").Value = 4844.5
$ws.Range("K107").Value = 4844.5
$ws.Range("M107").Value = -2924.5
$ws.Range("H134").Value = 2786.2727
$ws.Range("I134").Value = 2923.8333
$ws.Range("K134").Value = 8771.499899999999
$ws.Range("M134").Value = -6236.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28576408
$ws.Range("I31").Value = 4062.28
$ws.Range("K31").Value = 4062.28
$ws.Range("M31").Value = -3767.28
$ws.Range("H34").Value = 28576408
$ws.Range("I34").Value = 4062.28
$ws.Range("K34").Value = 4062.28
$ws.Range("M34").Value = -3860.28
$ws.Range("H43").Value = 45153.332
$ws.Range("J43").Value = 45153.332
$ws.Range("L43").Value = 45153.332
$ws.Range("N43").Value = -45521.332
$ws.Range("H101").Value = 45153.332
$ws.Range("J101").Value = 45153.332
$ws.Range("L101").Value = 45153.332
$ws.Range("N101").Value = -51643.332
$ws.Range("H122").Value = 2523.1667
$ws.Range("I122").Value = 2275.7144
$ws.Range("J122").Value = 2869.6
$ws.Range("K122").Value = 6827.1432
$ws.Range("L122").Value = 8608.799999999999
$ws.Range("M122").Value = -4377.1432
$ws.Range("N122").Value = -13508.8
$ws.Range("H132").Value = 71750.75999999999
$ws.Range("I132").Value = 93000.59
$ws.Range("K132").Value = 279001.77
$ws.Range("M132").Value = -276471.77
$ws.Range("H141").Value = 351068.8
$ws.Range("J141").Value = 365510.62
$ws.Range("L141").Value = 365510.62
$ws.Range("N141").Value = -375870.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3730
$ws.Range("J39").Value = 4940
$ws.Range("L39").Value = 14820
$ws.Range("N39").Value = -15408
$ws.Range("H55").Value = 4764550.5
$ws.Range("J55").Value = 3003.611
$ws.Range("L55").Value = 9010.832999999999
$ws.Range("N55").Value = -9364.832999999999
$ws.Range("H131").Value = 46660.04
$ws.Range("J131").Value = 7736.067
$ws.Range("L131").Value = 23208.201
$ws.Range("N131").Value = -33288.201
$ws.Range("H133").Value = 14562.833
$ws.Range("I133").Value = 9176
$ws.Range("J133").Value = 19949.666
$ws.Range("K133").Value = 27528
$ws.Range("L133").Value = 59848.99800000001
$ws.Range("M133").Value = -22468
$ws.Range("N133").Value = -69968.99800000001
$ws.Range("H134").Value = 3655.6667
$ws.Range("I134").Value = 1675.1875
$ws.Range("K134").Value = 5025.5625
$ws.Range("M134").Value = 44.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10280
$ws.Range("H80").Value = 5080.727
$ws.Range("I80").Value = 4735
$ws.Range("J80").Value = 5210.375
$ws.Range("K80").Value = 4735
$ws.Range("L80").Value = 5210.375
$ws.Range("M80").Value = -3737
$ws.Range("N80").Value = -7206.375
$ws.Range("H83").Value = 5080.727
$ws.Range("I83").Value = 4735
$ws.Range("J83").Value = 5210.375
$ws.Range("K83").Value = 23675
$ws.Range("L83").Value = 26051.875
$ws.Range("M83").Value = -18683
$ws.Range("N83").Value = -36035.875
$ws.Range("H97").Value = 2418.4614
$ws.Range("I97").Value = 1194.1
$ws.Range("K97").Value = 1194.1
$ws.Range("M97").Value = -698.0999999999999
$ws.Range("H102").Value = 1680.5186
$ws.Range("I102").Value = 1516.44
$ws.Range("J102").Value = 3731.5
$ws.Range("K102").Value = 1516.44
$ws.Range("L102").Value = 3731.5
$ws.Range("M102").Value = 105.5599999999999
$ws.Range("N102").Value = -6975.5
$ws.Range("H113").Value = 6484.4443
$ws.Range("I113").Value = 7999.3335
$ws.Range("J113").Value = 5727
$ws.Range("K113").Value = 7999.3335
$ws.Range("L113").Value = 5727
$ws.Range("M113").Value = -5829.3335
$ws.Range("N113").Value = -10067
$ws.Range("H126").Value = 5361658.5
$ws.Range("I126").Value = 2504285.8
$ws.Range("J126").Value = 12505090
$ws.Range("K126").Value = 7512857.399999999
$ws.Range("L126").Value = 37515270
$ws.Range("M126").Value = -7510387.399999999
$ws.Range("N126").Value = -37520210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4102.719
$ws.Range("I7").Value = 3619.8125
$ws.Range("J7").Value = 4720.84
$ws.Range("K7").Value = 3619.8125
$ws.Range("L7").Value = 4720.84
$ws.Range("M7").Value = -3507.8125
$ws.Range("N7").Value = -4944.84
$ws.Range("H40").Value = 5526.28
$ws.Range("I40").Value = 5436.619
$ws.Range("K40").Value = 5436.619
$ws.Range("M40").Value = -5300.619
$ws.Range("H46").Value = 1339.421
$ws.Range("J46").Value = 3716.3333
$ws.Range("L46").Value = 3716.3333
$ws.Range("N46").Value = -4092.3333
$ws.Range("H61").Value = 4749
$ws.Range("I61").Value = 2796.6667
$ws.Range("J61").Value = 6213.25
$ws.Range("K61").Value = 2796.6667
$ws.Range("L61").Value = 6213.25
$ws.Range("M61").Value = -2594.6667
$ws.Range("N61").Value = -6617.25
$ws.Range("H68").Value = 5385
$ws.Range("I68").Value = 3900
$ws.Range("K68").Value = 3900
$ws.Range("M68").Value = -3151
$ws.Range("H71").Value = 5385
$ws.Range("I71").Value = 3900
$ws.Range("K71").Value = 19500
$ws.Range("M71").Value = -15756
$ws.Range("H93").Value = 1768.1875
$ws.Range("I93").Value = 1306.2142
$ws.Range("J93").Value = 5002
$ws.Range("K93").Value = 1306.2142
$ws.Range("L93").Value = 5002
$ws.Range("M93").Value = -58.21419999999989
$ws.Range("N93").Value = -7498
$ws.Range("H113").Value = 4749
$ws.Range("I113").Value = 2796.6667
$ws.Range("J113").Value = 6213.25
$ws.Range("K113").Value = 2796.6667
$ws.Range("L113").Value = 6213.25
$ws.Range("M113").Value = -626.6667000000002
$ws.Range("N113").Value = -10553.25
$ws.Range("H122").Value = 4932
$ws.Range("J122").Value = 5833
$ws.Range("L122").Value = 17499
$ws.Range("N122").Value = -22399
$ws.Range("H126").Value = 4102.719
$ws.Range("I126").Value = 3619.8125
$ws.Range("J126").Value = 4720.84
$ws.Range("K126").Value = 10859.4375
$ws.Range("L126").Value = 14162.52
$ws.Range("M126").Value = -8389.4375
$ws.Range("N126").Value = -19102.52
$ws.Range("H132").Value = 38465252
$ws.Range("I132").Value = 3510.081
$ws.Range("K132").Value = 10530.243
$ws.Range("M132").Value = -8000.243
$ws.Range("H136").Value = 5290.696
$ws.Range("I136").Value = 4710
$ws.Range("J136").Value = 9162
$ws.Range("K136").Value = 14130
$ws.Range("L136").Value = 27486
$ws.Range("M136").Value = -11580
$ws.Range("N136").Value = -32586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 150
$ws.Range("I11").Value = 150
$ws.Range("K11").Value = 150
$ws.Range("M11").Value = -8
$ws.Range("H20").Value = 7333.3335
$ws.Range("I20").Value = 6000
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = -5760
$ws.Range("H52").Value = 27248.5
$ws.Range("I52").Value = 23749.25
$ws.Range("K52").Value = 23749.25
$ws.Range("M52").Value = -23523.25
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -26476
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H86").Value = 11178588
$ws.Range("I86").Value = 25036824
$ws.Range("J86").Value = 91999.8
$ws.Range("K86").Value = 25036824
$ws.Range("L86").Value = 91999.8
$ws.Range("M86").Value = -25035701
$ws.Range("N86").Value = -94245.8
$ws.Range("H89").Value = 11178588
$ws.Range("I89").Value = 25036824
$ws.Range("J89").Value = 91999.8
$ws.Range("K89").Value = 125184120
$ws.Range("L89").Value = 459999
$ws.Range("M89").Value = -125178504
$ws.Range("N89").Value = -471231
$ws.Range("H96").Value = 6928.0586
$ws.Range("I96").Value = 6940.2
$ws.Range("J96").Value = 6923
$ws.Range("K96").Value = 6940.2
$ws.Range("L96").Value = 6923
$ws.Range("M96").Value = -5567.2
$ws.Range("N96").Value = -9669
$ws.Range("H126").Value = 5936.3335
$ws.Range("I126").Value = 6539.5293
$ws.Range("J126").Value = 4471.4287
$ws.Range("K126").Value = 19618.5879
$ws.Range("L126").Value = 13414.2861
$ws.Range("M126").Value = -17148.5879
$ws.Range("N126").Value = -18354.2861
$ws.Range("H132").Value = 4349.095
$ws.Range("I132").Value = 4911.2827
$ws.Range("J132").Value = 2827.8823
$ws.Range("K132").Value = 14733.8481
$ws.Range("L132").Value = 8483.6469
$ws.Range("M132").Value = -12203.8481
$ws.Range("N132").Value = -13543.6469

Write-Output "Applied all changes"